# Add population for accIsActive[Account*Account]
# New relation column J on the "Identity Provider data" sheet:
#   J12 = header text "accIsActive"
#   J13 = type text "Account"
#   J14:J18 = formula referencing back to the row's own Account name (A<row>)
# Also: row 14 (Stef) becomes an AccountMgr too (H14 = "AccountMgr"),
# matching rows 16/18 which already carry that allowed role.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / type row labels for the new column.
$ws.Range("J12").Value = "accIsActive"
$ws.Range("J13").Value = "Account"

# Give the new data cells the same centered style used by the sibling
# G/H/I columns before filling them in, so the xf matches (style index 3).
$ws.Range("H14").HorizontalAlignment = -4108
$ws.Range("J14:J18").HorizontalAlignment = -4108

# Stef's account also gets the AccountMgr allowed role.
$ws.Range("H14").Value = "AccountMgr"

# Fill the new self-referencing formula column.
$ws.Range("J14").Formula = "=A14"
$ws.Range("J15:J18").Formula = "=A15"
